# Refresh admin data loader output: update population/average/percentile series
# for the 2016, 2017, 2018 sheets of total-PER.xlsx (setup admin data loaders).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("2016")
$ws.Range("C2").Value = 30866494
$ws.Range("D2").Value = 5515.40966796875
$ws.Range("E2").Value = 0.75727224349975586
$ws.Range("C3").Value = 30866494
$ws.Range("D3").Value = 5515.40966796875
$ws.Range("E3").Value = 0.8230818510055542
$ws.Range("C4").Value = 30866494
$ws.Range("D4").Value = 5515.40966796875
$ws.Range("E4").Value = 0.8585284948348999
$ws.Range("C5").Value = 30866494
$ws.Range("D5").Value = 5515.40966796875
$ws.Range("E5").Value = 0.92302888631820679
$ws.Range("C6").Value = 30866494
$ws.Range("D6").Value = 5515.40966796875
$ws.Range("E6").Value = 0.95472633838653564
$ws.Range("C7").Value = 30866494
$ws.Range("D7").Value = 5515.40966796875
$ws.Range("E7").Value = 0.97717702388763428
$ws.Range("C8").Value = 30866494
$ws.Range("D8").Value = 5515.40966796875
$ws.Range("E8").Value = 0.98544025421142578
$ws.Range("C9").Value = 30866494
$ws.Range("D9").Value = 5515.40966796875
$ws.Range("E9").Value = 0.99143153429031372
$ws.Range("C10").Value = 30866494
$ws.Range("D10").Value = 5515.40966796875
$ws.Range("E10").Value = 0.99571722745895386
$ws.Range("C11").Value = 30866494
$ws.Range("D11").Value = 5515.40966796875
$ws.Range("E11").Value = 0.99748361110687256
$ws.Range("C12").Value = 30866494
$ws.Range("D12").Value = 5515.40966796875
$ws.Range("E12").Value = 0.9988057017326355
$ws.Range("C13").Value = 30866494
$ws.Range("D13").Value = 5515.40966796875
$ws.Range("E13").Value = 0.99926990270614624
$ws.Range("C14").Value = 30866494
$ws.Range("D14").Value = 5515.40966796875
$ws.Range("E14").Value = 0.99962466955184937
$ws.Range("C15").Value = 30866494
$ws.Range("D15").Value = 5515.40966796875
$ws.Range("E15").Value = 0.99976241588592529
$ws.Range("C16").Value = 30866494
$ws.Range("D16").Value = 5515.40966796875
$ws.Range("E16").Value = 0.99983477592468262
$ws.Range("C17").Value = 30866494
$ws.Range("D17").Value = 5515.40966796875
$ws.Range("E17").Value = 0.9999430775642395
$ws.Range("C18").Value = 30866494
$ws.Range("D18").Value = 5515.40966796875
$ws.Range("E18").Value = 0.9999699592590332

$ws = $wb.Worksheets.Item("2017")
$ws.Range("C2").Value = 31324636
$ws.Range("D2").Value = 5971.404296875
$ws.Range("E2").Value = 0.7574123740196228
$ws.Range("C3").Value = 31324636
$ws.Range("D3").Value = 5971.404296875
$ws.Range("E3").Value = 0.8210064172744751
$ws.Range("C4").Value = 31324636
$ws.Range("D4").Value = 5971.404296875
$ws.Range("E4").Value = 0.85172677040100098
$ws.Range("C5").Value = 31324636
$ws.Range("D5").Value = 5971.404296875
$ws.Range("E5").Value = 0.91912400722503662
$ws.Range("C6").Value = 31324636
$ws.Range("D6").Value = 5971.404296875
$ws.Range("E6").Value = 0.95218104124069214
$ws.Range("C7").Value = 31324636
$ws.Range("D7").Value = 5971.404296875
$ws.Range("E7").Value = 0.97609686851501465
$ws.Range("C8").Value = 31324636
$ws.Range("D8").Value = 5971.404296875
$ws.Range("E8").Value = 0.98486131429672241
$ws.Range("C9").Value = 31324636
$ws.Range("D9").Value = 5971.404296875
$ws.Range("E9").Value = 0.99109137058258057
$ws.Range("C10").Value = 31324636
$ws.Range("D10").Value = 5971.404296875
$ws.Range("E10").Value = 0.99559754133224487
$ws.Range("C11").Value = 31324636
$ws.Range("D11").Value = 5971.404296875
$ws.Range("E11").Value = 0.99740010499954224
$ws.Range("C12").Value = 31324636
$ws.Range("D12").Value = 5971.404296875
$ws.Range("E12").Value = 0.99875748157501221
$ws.Range("C13").Value = 31324636
$ws.Range("D13").Value = 5971.404296875
$ws.Range("E13").Value = 0.99924463033676147
$ws.Range("C14").Value = 31324636
$ws.Range("D14").Value = 5971.404296875
$ws.Range("E14").Value = 0.99962013959884644
$ws.Range("C15").Value = 31324636
$ws.Range("D15").Value = 5971.404296875
$ws.Range("E15").Value = 0.99976116418838501
$ws.Range("C16").Value = 31324636
$ws.Range("D16").Value = 5971.404296875
$ws.Range("E16").Value = 0.99983441829681396
$ws.Range("C17").Value = 31324636
$ws.Range("D17").Value = 5971.404296875
$ws.Range("E17").Value = 0.99994409084320068
$ws.Range("C18").Value = 31324636
$ws.Range("D18").Value = 5971.404296875
$ws.Range("E18").Value = 0.99996960163116455

$ws = $wb.Worksheets.Item("2018")
$ws.Range("C2").Value = 31897584
$ws.Range("D2").Value = 6542.81787109375
$ws.Range("E2").Value = 0.75383388996124268
$ws.Range("C3").Value = 31897584
$ws.Range("D3").Value = 6542.81787109375
$ws.Range("E3").Value = 0.81654459238052368
$ws.Range("C4").Value = 31897584
$ws.Range("D4").Value = 6542.81787109375
$ws.Range("E4").Value = 0.84587550163269043
$ws.Range("C5").Value = 31897584
$ws.Range("D5").Value = 6542.81787109375
$ws.Range("E5").Value = 0.91002506017684937
$ws.Range("C6").Value = 31897584
$ws.Range("D6").Value = 6542.81787109375
$ws.Range("E6").Value = 0.9477723240852356
$ws.Range("C7").Value = 31897584
$ws.Range("D7").Value = 6542.81787109375
$ws.Range("E7").Value = 0.97422713041305542
$ws.Range("C8").Value = 31897584
$ws.Range("D8").Value = 6542.81787109375
$ws.Range("E8").Value = 0.98377478122711182
$ws.Range("C9").Value = 31897584
$ws.Range("D9").Value = 6542.81787109375
$ws.Range("E9").Value = 0.99035501480102539
$ws.Range("C10").Value = 31897584
$ws.Range("D10").Value = 6542.81787109375
$ws.Range("E10").Value = 0.99517327547073364
$ws.Range("C11").Value = 31897584
$ws.Range("D11").Value = 6542.81787109375
$ws.Range("E11").Value = 0.99711304903030396
$ws.Range("C12").Value = 31897584
$ws.Range("D12").Value = 6542.81787109375
$ws.Range("E12").Value = 0.99864649772644043
$ws.Range("C13").Value = 31897584
$ws.Range("D13").Value = 6542.81787109375
$ws.Range("E13").Value = 0.99918973445892334
$ws.Range("C14").Value = 31897584
$ws.Range("D14").Value = 6542.81787109375
$ws.Range("E14").Value = 0.99959456920623779
$ws.Range("C15").Value = 31897584
$ws.Range("D15").Value = 6542.81787109375
$ws.Range("E15").Value = 0.99974685907363892
$ws.Range("C16").Value = 31897584
$ws.Range("D16").Value = 6542.81787109375
$ws.Range("E16").Value = 0.99982273578643799
$ws.Range("C17").Value = 31897584
$ws.Range("D17").Value = 6542.81787109375
$ws.Range("E17").Value = 0.9999358057975769
$ws.Range("C18").Value = 31897584
$ws.Range("D18").Value = 6542.81787109375
$ws.Range("E18").Value = 0.99996423721313477

